$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 881
$ws.Cells.Item(2, 10).Value = 1915.6666
$ws.Cells.Item(2, 12).Value = 1915.6666
$ws.Cells.Item(2, 14).Value = -2141.6666

$ws.Cells.Item(15, 8).Value = 1096.6232
$ws.Cells.Item(15, 9).Value = 1096.6232
$ws.Cells.Item(15, 11).Value = 3289.8696
$ws.Cells.Item(15, 13).Value = -3120.8696

$ws.Cells.Item(40, 8).Value = 83335540
$ws.Cells.Item(40, 9).Value = 2412.6667
$ws.Cells.Item(40, 11).Value = 2412.6667
$ws.Cells.Item(40, 13).Value = -2237.6667

$ws.Cells.Item(74, 8).Value = 6961
$ws.Cells.Item(74, 9).Value = 2993
$ws.Cells.Item(74, 11).Value = 2993
$ws.Cells.Item(74, 13).Value = -2057

$ws.Cells.Item(77, 8).Value = 6961
$ws.Cells.Item(77, 9).Value = 2993
$ws.Cells.Item(77, 11).Value = 14965
$ws.Cells.Item(77, 13).Value = -10285

$ws.Cells.Item(99, 8).Value = 4010.9092
$ws.Cells.Item(99, 9).Value = 324.6
$ws.Cells.Item(99, 10).Value = 7082.8335
$ws.Cells.Item(99, 11).Value = 973.8000000000001
$ws.Cells.Item(99, 12).Value = 21248.5005
$ws.Cells.Item(99, 13).Value = 524.1999999999999
$ws.Cells.Item(99, 14).Value = -24244.5005

$ws.Cells.Item(100, 8).Value = 5479.8335
$ws.Cells.Item(100, 9).Value = 4606.1665
$ws.Cells.Item(100, 11).Value = 4606.1665
$ws.Cells.Item(100, 13).Value = -4065.1665

$ws.Cells.Item(101, 8).Value = 946.63635
$ws.Cells.Item(101, 10).Value = 1550
$ws.Cells.Item(101, 12).Value = 4650
$ws.Cells.Item(101, 14).Value = -7894

$ws.Cells.Item(112, 8).Value = 5805.5835
$ws.Cells.Item(112, 10).Value = 6242.4546
$ws.Cells.Item(112, 12).Value = 18727.3638
$ws.Cells.Item(112, 14).Value = -20943.3638

$ws.Cells.Item(116, 8).Value = 11144.852
$ws.Cells.Item(116, 10).Value = 12798.4
$ws.Cells.Item(116, 12).Value = 12798.4
$ws.Cells.Item(116, 14).Value = -19682.4

$ws.Cells.Item(138, 8).Value = 3313.6
$ws.Cells.Item(138, 9).Value = 2840.1428
$ws.Cells.Item(138, 10).Value = 3390.6743
$ws.Cells.Item(138, 11).Value = 8520.428400000001
$ws.Cells.Item(138, 12).Value = 10172.0229
$ws.Cells.Item(138, 13).Value = -3380.428400000001
$ws.Cells.Item(138, 14).Value = -20452.0229

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9447.734
$ws.Cells.Item(32, 9).Value = 9049.308000000001
$ws.Cells.Item(32, 11).Value = 9049.308000000001
$ws.Cells.Item(32, 13).Value = -8762.308000000001

$ws.Cells.Item(53, 8).Value = 19930
$ws.Cells.Item(53, 9).Value = 9895.5
$ws.Cells.Item(53, 11).Value = 9895.5
$ws.Cells.Item(53, 13).Value = -9213.5

$ws.Cells.Item(61, 8).Value = 5389101
$ws.Cells.Item(61, 9).Value = 7147755
$ws.Cells.Item(61, 10).Value = 912526.4
$ws.Cells.Item(61, 11).Value = 7147755
$ws.Cells.Item(61, 12).Value = 912526.4
$ws.Cells.Item(61, 13).Value = -7147543
$ws.Cells.Item(61, 14).Value = -912950.4

$ws.Cells.Item(97, 8).Value = 1814.6451
$ws.Cells.Item(97, 10).Value = 5301.5713
$ws.Cells.Item(97, 12).Value = 5301.5713
$ws.Cells.Item(97, 14).Value = -6293.5713

$ws.Cells.Item(102, 8).Value = 2559
$ws.Cells.Item(102, 9).Value = 1830.7059
$ws.Cells.Item(102, 11).Value = 1830.7059
$ws.Cells.Item(102, 13).Value = -208.7058999999999

$ws.Cells.Item(110, 8).Value = 6006.1816
$ws.Cells.Item(110, 9).Value = 6346.6665
$ws.Cells.Item(110, 10).Value = 5276.5713
$ws.Cells.Item(110, 11).Value = 6346.6665
$ws.Cells.Item(110, 12).Value = 5276.5713
$ws.Cells.Item(110, 13).Value = -4301.6665
$ws.Cells.Item(110, 14).Value = -9366.5713

$ws.Cells.Item(136, 8).Value = 5389101
$ws.Cells.Item(136, 9).Value = 7147755
$ws.Cells.Item(136, 10).Value = 912526.4
$ws.Cells.Item(136, 11).Value = 21443265
$ws.Cells.Item(136, 12).Value = 2737579.2
$ws.Cells.Item(136, 13).Value = -21440715
$ws.Cells.Item(136, 14).Value = -2742679.2

$ws.Cells.Item(139, 8).Value = 148975.75
$ws.Cells.Item(139, 10).Value = 157973.55
$ws.Cells.Item(139, 12).Value = 157973.55
$ws.Cells.Item(139, 14).Value = -168253.55

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 2110.0303
$ws.Cells.Item(94, 9).Value = 2377.0417
$ws.Cells.Item(94, 10).Value = 1398
$ws.Cells.Item(94, 11).Value = 2377.0417
$ws.Cells.Item(94, 12).Value = 1398
$ws.Cells.Item(94, 13).Value = -1926.0417
$ws.Cells.Item(94, 14).Value = -2300

$ws.Cells.Item(134, 8).Value = 3335591.5
$ws.Cells.Item(134, 9).Value = 1679
$ws.Cells.Item(134, 11).Value = 5037
$ws.Cells.Item(134, 13).Value = -2502

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 22314.7
$ws.Cells.Item(62, 9).Value = 14749.4
$ws.Cells.Item(62, 10).Value = 29880
$ws.Cells.Item(62, 11).Value = 14749.4
$ws.Cells.Item(62, 12).Value = 29880
$ws.Cells.Item(62, 13).Value = -14125.4
$ws.Cells.Item(62, 14).Value = -31128

$ws.Cells.Item(65, 8).Value = 22314.7
$ws.Cells.Item(65, 9).Value = 14749.4
$ws.Cells.Item(65, 10).Value = 29880
$ws.Cells.Item(65, 11).Value = 73747
$ws.Cells.Item(65, 12).Value = 149400
$ws.Cells.Item(65, 13).Value = -70627
$ws.Cells.Item(65, 14).Value = -155640

$ws.Cells.Item(100, 8).Value = 97559.60000000001
$ws.Cells.Item(100, 10).Value = 97559.60000000001
$ws.Cells.Item(100, 12).Value = 97559.60000000001
$ws.Cells.Item(100, 14).Value = -99723.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 122.29412
$ws.Cells.Item(2, 9).Value = 114.333336
$ws.Cells.Item(2, 10).Value = 141.4
$ws.Cells.Item(2, 11).Value = 686.000016
$ws.Cells.Item(2, 12).Value = 848.4000000000001
$ws.Cells.Item(2, 13).Value = -573.000016
$ws.Cells.Item(2, 14).Value = -1074.4

$ws.Cells.Item(23, 8).Value = 368.13333
$ws.Cells.Item(23, 9).Value = 192.2
$ws.Cells.Item(23, 10).Value = 456.1
$ws.Cells.Item(23, 11).Value = 576.5999999999999
$ws.Cells.Item(23, 12).Value = 1368.3
$ws.Cells.Item(23, 13).Value = -341.5999999999999
$ws.Cells.Item(23, 14).Value = -1838.3

$ws.Cells.Item(38, 8).Value = 30.833334
$ws.Cells.Item(38, 10).Value = 54
$ws.Cells.Item(38, 12).Value = 162
$ws.Cells.Item(38, 14).Value = -856

$ws.Cells.Item(131, 8).Value = 3128.6094
$ws.Cells.Item(131, 9).Value = 2674.4
$ws.Cells.Item(131, 10).Value = 3267.653
$ws.Cells.Item(131, 11).Value = 8023.200000000001
$ws.Cells.Item(131, 12).Value = 9802.958999999999
$ws.Cells.Item(131, 13).Value = -2983.200000000001
$ws.Cells.Item(131, 14).Value = -19882.959

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(103, 8).Value = 85000
$ws.Cells.Item(103, 10).Value = 85000
$ws.Cells.Item(103, 12).Value = 85000
$ws.Cells.Item(103, 14).Value = -87344

$ws.Cells.Item(123, 8).Value = 45320
$ws.Cells.Item(123, 10).Value = 45320
$ws.Cells.Item(123, 12).Value = 45320
$ws.Cells.Item(123, 14).Value = -50220

$ws.Cells.Item(132, 8).Value = 2383844.8
$ws.Cells.Item(132, 9).Value = 2864.3513
$ws.Cells.Item(132, 10).Value = 20003098
$ws.Cells.Item(132, 11).Value = 8593.053899999999
$ws.Cells.Item(132, 12).Value = 60009294
$ws.Cells.Item(132, 13).Value = -6063.053899999999
$ws.Cells.Item(132, 14).Value = -60014354

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 26401950
$ws.Cells.Item(22, 9).Value = 66000324
$ws.Cells.Item(22, 11).Value = 66000324
$ws.Cells.Item(22, 13).Value = -66000029

$ws.Cells.Item(27, 8).Value = 26401950
$ws.Cells.Item(27, 9).Value = 66000324
$ws.Cells.Item(27, 11).Value = 66000324
$ws.Cells.Item(27, 13).Value = -66000217

$ws.Cells.Item(57, 8).Value = 28041.143
$ws.Cells.Item(57, 9).Value = 28041.143
$ws.Cells.Item(57, 11).Value = 28041.143
$ws.Cells.Item(57, 13).Value = -27475.143
